$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2025-07-27 Sunday", $false, $false, $false, $false, $false, $true, 1, $false, "2025-07-28 Monday", 2) | Out-Null

# Update each math expression cell in the table (20 rows x 5 columns)
$tbl = $d.Tables.Item(1)

$tbl.Cell(1, 1).Range.Text = "40+17="  # was "66-3="
$tbl.Cell(1, 2).Range.Text = "55-48="  # was "82-57="
$tbl.Cell(1, 3).Range.Text = "93-29="  # was "46-1="
$tbl.Cell(1, 4).Range.Text = "9+83="  # was "47-40="
$tbl.Cell(1, 5).Range.Text = "35+48="  # was "68+25="

$tbl.Cell(2, 1).Range.Text = "31+24="  # was "72-35="
$tbl.Cell(2, 2).Range.Text = "1+54="  # was "83-27="
$tbl.Cell(2, 3).Range.Text = "45-34="  # was "5+21="
$tbl.Cell(2, 4).Range.Text = "46+43="  # was "98-85="
$tbl.Cell(2, 5).Range.Text = "81+16="  # was "76-4="

$tbl.Cell(3, 1).Range.Text = "13+80="  # was "53+19="
$tbl.Cell(3, 2).Range.Text = "65+17="  # was "79-2="
$tbl.Cell(3, 3).Range.Text = "38+12="  # was "76-33="
$tbl.Cell(3, 4).Range.Text = "51-19="  # was "74-50="
$tbl.Cell(3, 5).Range.Text = "5+80="  # was "57-22="

$tbl.Cell(4, 1).Range.Text = "57-27="  # was "23+22="
$tbl.Cell(4, 2).Range.Text = "89-11="  # was "1+24="
$tbl.Cell(4, 3).Range.Text = "40+51="  # was "67-65="
$tbl.Cell(4, 4).Range.Text = "15+71="  # was "41+3="
$tbl.Cell(4, 5).Range.Text = "38+3="  # was "9+34="

$tbl.Cell(5, 1).Range.Text = "61-54="  # was "73-18="
$tbl.Cell(5, 2).Range.Text = "51+10="  # was "88-17="
$tbl.Cell(5, 3).Range.Text = "97-53="  # was "34+13="
$tbl.Cell(5, 4).Range.Text = "53-2="  # was "85-59="
$tbl.Cell(5, 5).Range.Text = "16-15="  # was "29-3="

$tbl.Cell(6, 1).Range.Text = "67-60="  # was "89-53="
$tbl.Cell(6, 2).Range.Text = "70+29="  # was "18+70="
$tbl.Cell(6, 3).Range.Text = "4+7="  # was "24+73="
$tbl.Cell(6, 4).Range.Text = "33+56="  # was "30-18="
$tbl.Cell(6, 5).Range.Text = "90-88="  # was "0+40="

$tbl.Cell(7, 1).Range.Text = "12+14="  # was "18-0="
$tbl.Cell(7, 2).Range.Text = "68-26="  # was "77+9="
$tbl.Cell(7, 3).Range.Text = "97-84="  # was "2+42="
$tbl.Cell(7, 4).Range.Text = "20+65="  # was "43+3="
$tbl.Cell(7, 5).Range.Text = "98-81="  # was "96-70="

$tbl.Cell(8, 1).Range.Text = "3+47="  # was "74-67="
$tbl.Cell(8, 2).Range.Text = "45+33="  # was "53-39="
$tbl.Cell(8, 3).Range.Text = "85-71="  # was "24+9="
$tbl.Cell(8, 4).Range.Text = "39+5="  # was "11+6="
$tbl.Cell(8, 5).Range.Text = "49+1="  # was "22+15="

$tbl.Cell(9, 1).Range.Text = "38+1="  # was "57-37="
$tbl.Cell(9, 2).Range.Text = "59-57="  # was "95-45="
$tbl.Cell(9, 3).Range.Text = "32+35="  # was "34+36="
$tbl.Cell(9, 4).Range.Text = "93-80="  # was "74-41="
$tbl.Cell(9, 5).Range.Text = "73-37="  # was "83-40="

$tbl.Cell(10, 1).Range.Text = "62+15="  # was "89-41="
$tbl.Cell(10, 2).Range.Text = "24+28="  # was "71-68="
$tbl.Cell(10, 3).Range.Text = "96-88="  # was "39+47="
$tbl.Cell(10, 4).Range.Text = "24+61="  # was "4+75="
$tbl.Cell(10, 5).Range.Text = "84-47="  # was "92-8="

$tbl.Cell(11, 1).Range.Text = "97-2="  # was "91-57="
$tbl.Cell(11, 2).Range.Text = "24+1="  # was "32+45="
$tbl.Cell(11, 3).Range.Text = "64-24="  # was "49-19="
$tbl.Cell(11, 4).Range.Text = "46+47="  # was "61+36="
$tbl.Cell(11, 5).Range.Text = "11+18="  # was "56-2="

$tbl.Cell(12, 1).Range.Text = "41-35="  # was "36+49="
$tbl.Cell(12, 2).Range.Text = "4+48="  # was "11-2="
$tbl.Cell(12, 3).Range.Text = "59-54="  # was "95-80="
$tbl.Cell(12, 4).Range.Text = "41-0="  # was "91-25="
$tbl.Cell(12, 5).Range.Text = "97-68="  # was "65-16="

$tbl.Cell(13, 1).Range.Text = "89-22="  # was "85-20="
$tbl.Cell(13, 2).Range.Text = "22+34="  # was "48-42="
$tbl.Cell(13, 3).Range.Text = "31+8="  # was "69-11="
$tbl.Cell(13, 4).Range.Text = "62-27="  # was "23+41="
$tbl.Cell(13, 5).Range.Text = "19+36="  # was "93-47="

$tbl.Cell(14, 1).Range.Text = "99-10="  # was "70-23="
$tbl.Cell(14, 2).Range.Text = "41+44="  # was "50-3="
$tbl.Cell(14, 3).Range.Text = "78-63="  # was "25+73="
$tbl.Cell(14, 4).Range.Text = "83-62="  # was "54+10="
$tbl.Cell(14, 5).Range.Text = "39-32="  # was "92-13="

$tbl.Cell(15, 1).Range.Text = "43+19="  # was "11+69="
$tbl.Cell(15, 2).Range.Text = "3+21="  # was "90-69="
$tbl.Cell(15, 3).Range.Text = "46-22="  # was "49+12="
$tbl.Cell(15, 4).Range.Text = "29-16="  # was "2+46="
$tbl.Cell(15, 5).Range.Text = "34-21="  # was "89-64="

$tbl.Cell(16, 1).Range.Text = "86-62="  # was "87-82="
$tbl.Cell(16, 2).Range.Text = "78+10="  # was "39+60="
$tbl.Cell(16, 3).Range.Text = "63-31="  # was "39+42="
$tbl.Cell(16, 4).Range.Text = "41+38="  # was "19+4="
$tbl.Cell(16, 5).Range.Text = "94+1="  # was "91-26="

$tbl.Cell(17, 1).Range.Text = "79+19="  # was "19-8="
$tbl.Cell(17, 2).Range.Text = "14+15="  # was "80-71="
$tbl.Cell(17, 3).Range.Text = "45-40="  # was "5+11="
$tbl.Cell(17, 4).Range.Text = "86-14="  # was "76-21="
$tbl.Cell(17, 5).Range.Text = "80-73="  # was "40-36="

$tbl.Cell(18, 1).Range.Text = "44+35="  # was "32+43="
$tbl.Cell(18, 2).Range.Text = "7+36="  # was "86-27="
$tbl.Cell(18, 3).Range.Text = "91-61="  # was "16+52="
$tbl.Cell(18, 4).Range.Text = "96-32="  # was "68-46="
$tbl.Cell(18, 5).Range.Text = "59+16="  # was "94-30="

$tbl.Cell(19, 1).Range.Text = "99-66="  # was "17-2="
$tbl.Cell(19, 2).Range.Text = "63-51="  # was "17+41="
$tbl.Cell(19, 3).Range.Text = "6+75="  # was "16+62="
$tbl.Cell(19, 4).Range.Text = "19+61="  # was "0+17="
$tbl.Cell(19, 5).Range.Text = "92-17="  # was "93-35="

$tbl.Cell(20, 1).Range.Text = "24-20="  # was "47-35="
$tbl.Cell(20, 2).Range.Text = "53-20="  # was "56-4="
$tbl.Cell(20, 3).Range.Text = "39+40="  # was "96-60="
$tbl.Cell(20, 4).Range.Text = "40+23="  # was "8+66="
$tbl.Cell(20, 5).Range.Text = "0+52="  # was "53-41="
